$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename destination barcode suffix from _C to _E (rows 2-7)
$ws.Range("E2").Value = "ssdest000000141jul17_E"
$ws.Range("E3").Value = "ssdest000000141jul17_E"
$ws.Range("E4").Value = "ssdest000000141jul17_E"
$ws.Range("E5").Value = "ssdest000000141jul17_E"
$ws.Range("E6").Value = "ssdest000000141jul17_384_E"
$ws.Range("E7").Value = "ssdest000000141jul17_384_E"

# Row 6's source well id changes from A7 to A6
$ws.Range("C6").Value = "A6"

# Add new "Dest Well Count" column H
$ws.Columns("H").ColumnWidth = 25.5

# Give the new header cell the same look as the rest of the header row
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Dest Well Count"
$ws.Range("H2").Value = 96
$ws.Range("H3").Value = 96
$ws.Range("H4").Value = 96
$ws.Range("H5").Value = 96
$ws.Range("H6").Value = 384
$ws.Range("H7").Value = 384

# Match the selection/active cell shown after the edit
$ws.Range("H7").Select()
